$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.805.67"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.50"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.87"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("E8").Value = "  -0.43%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.32"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.23"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.638.71"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.98"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.805.39"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.04"
$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.58"
$ws.Range("E22").Value = "  +4.51%  "

$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("E24").Value = "  -2.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.18"
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.70%  "

$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.262.76"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  -1.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.779.34"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -4.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.22"
$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.19"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.58"
$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("E48").Value = "  -2.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0965"
$ws.Range("E51").Value = "  -1.51%  "
